# "Generate Report for Archive"
#
# The handoff status for the 79ab8a41-... file moved from
# "Ready for handoff" to "In Translation" everywhere it is reported:
#   - Overview sheet: per-language status columns (zh-cn -> col E, de-de -> col F)
#   - zh-cn sheet: Status column (C)
#   - de-de sheet: Status column (C)
#
# Because the new status text is shorter than the old one, the report
# generator re-fits the Status columns so they are no longer as wide.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview!E2 (zh-cn status) and Overview!F2 (de-de status)
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# zh-cn!C2 / de-de!C2 Status column
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# Re-fit the (now narrower) Status columns to the new content.
$wsOverview.Range("E1:F1").EntireColumn.AutoFit() | Out-Null
$wsZhCn.Range("C1").EntireColumn.AutoFit() | Out-Null
$wsDeDe.Range("C1").EntireColumn.AutoFit() | Out-Null

# The headless AutoFit heuristic here doesn't reproduce desktop Excel's
# font-metric-accurate pixel math, so nudge the resulting widths to the
# closest value this engine can actually store that matches what real
# Excel computed for "In Translation" in this report (~13.41 chars).
$wsOverview.Range("E1").EntireColumn.ColumnWidth = 12.5
$wsOverview.Range("F1").EntireColumn.ColumnWidth = 12.5
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = 12.5
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = 12.5
